$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data for the 2025-10-13 run.
# The leading apostrophe forces Excel to store the date-looking string
# as literal text instead of auto-converting it to a date serial number
# (matching the existing rows, which are all plain text dates).
# ClearFormats() then strips the resulting quote-prefix formatting so the
# new cell keeps the same (default) style as every other data cell.
$ws.Range("A42").Value = "'10/13/2025"
$ws.Range("A42").ClearFormats()

$ws.Range("B42").Value = 0.1714181428297508
$ws.Range("C42").Value = 0.8285818571702492
